# Auto-generated Excel COM-interop edit script
# Applies numeric value updates (leve profit recalculations) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 40: Stuck in the Moment
$ws.Range("H40").Value = 1050.3334
$ws.Range("J40").Value = 1050.3334
$ws.Range("L40").Value = 1050.3334
$ws.Range("N40").Value = -1400.3334

# ALC row 107: Another Man's Ink
$ws.Range("H107").Value = 888.1111
$ws.Range("I107").Value = 1014
$ws.Range("J107").Value = 705
$ws.Range("K107").Value = 1014
$ws.Range("L107").Value = 705
$ws.Range("M107").Value = 906
$ws.Range("N107").Value = -4545

# ALC row 129: Practical Command
$ws.Range("H129").Value = 861.4
$ws.Range("J129").Value = 904.7646999999999
$ws.Range("L129").Value = 2714.2941
$ws.Range("N129").Value = -12714.2941

# ALC row 139: Something Salty and Ceremonial
$ws.Range("H139").Value = 74860
$ws.Range("J139").Value = 74860
$ws.Range("L139").Value = 74860
$ws.Range("N139").Value = -85140

$ws = $wb.Worksheets.Item("ARM")
# ARM row 2: Ain't Got No Ingots
$ws.Range("H2").Value = 641232.5600000001
$ws.Range("I2").Value = 1607.6786
$ws.Range("J2").Value = 1636204.6
$ws.Range("K2").Value = 1607.6786
$ws.Range("L2").Value = 1636204.6
$ws.Range("M2").Value = -1494.6786
$ws.Range("N2").Value = -1636430.6

# ARM row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 5559542
$ws.Range("I61").Value = 6948941
$ws.Range("J61").Value = 1945
$ws.Range("K61").Value = 6948941
$ws.Range("L61").Value = 1945
$ws.Range("M61").Value = -6948729
$ws.Range("N61").Value = -2369

# ARM row 63: Rivets Run through It
$ws.Range("H63").Value = 2545.28
$ws.Range("I63").Value = 1446.2222
$ws.Range("K63").Value = 1446.2222
$ws.Range("M63").Value = -760.2221999999999

# ARM row 66: A Riveting Revival (L)
$ws.Range("H66").Value = 2545.28
$ws.Range("I66").Value = 1446.2222
$ws.Range("K66").Value = 7231.111
$ws.Range("M66").Value = -3799.111

# ARM row 74: As the Bolt Flies
$ws.Range("H74").Value = 20002662
$ws.Range("I74").Value = 33334234
$ws.Range("J74").Value = 5304.1
$ws.Range("K74").Value = 33334234
$ws.Range("L74").Value = 5304.1
$ws.Range("M74").Value = -33333360
$ws.Range("N74").Value = -7052.1

# ARM row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 20002662
$ws.Range("I77").Value = 33334234
$ws.Range("J77").Value = 5304.1
$ws.Range("K77").Value = 166671170
$ws.Range("L77").Value = 26520.5
$ws.Range("M77").Value = -166666802
$ws.Range("N77").Value = -35256.5

# ARM row 116: No Scope
$ws.Range("H116").Value = 641232.5600000001
$ws.Range("I116").Value = 1607.6786
$ws.Range("J116").Value = 1636204.6
$ws.Range("K116").Value = 1607.6786
$ws.Range("L116").Value = 1636204.6
$ws.Range("M116").Value = 686.3214
$ws.Range("N116").Value = -1640792.6

# ARM row 136: Metal with Mettle
$ws.Range("H136").Value = 5559542
$ws.Range("I136").Value = 6948941
$ws.Range("J136").Value = 1945
$ws.Range("K136").Value = 20846823
$ws.Range("L136").Value = 5835
$ws.Range("M136").Value = -20844273
$ws.Range("N136").Value = -10935

$ws = $wb.Worksheets.Item("BSM")
# BSM row 3: Hells Bells
$ws.Range("H3").Value = 641232.5600000001
$ws.Range("I3").Value = 1607.6786
$ws.Range("J3").Value = 1636204.6
$ws.Range("K3").Value = 1607.6786
$ws.Range("L3").Value = 1636204.6
$ws.Range("M3").Value = -1493.6786
$ws.Range("N3").Value = -1636432.6

# BSM row 107: The Gold Experience
$ws.Range("H107").Value = 940873.4399999999
$ws.Range("I107").Value = 1007865.9
$ws.Range("J107").Value = 2980
$ws.Range("K107").Value = 1007865.9
$ws.Range("L107").Value = 2980
$ws.Range("M107").Value = -1005945.9
$ws.Range("N107").Value = -6820

# BSM row 134: Ruthenium Supremium
$ws.Range("H134").Value = 11117824
$ws.Range("I134").Value = 12507346
$ws.Range("J134").Value = 1650
$ws.Range("K134").Value = 37522038
$ws.Range("L134").Value = 4950
$ws.Range("M134").Value = -37519503
$ws.Range("N134").Value = -10020

# BSM row 140: Ceremonial Teeth
$ws.Range("H140").Value = 74780
$ws.Range("J140").Value = 74780
$ws.Range("L140").Value = 74780
$ws.Range("N140").Value = -85140

$ws = $wb.Worksheets.Item("CRP")
# CRP row 31: Wall Not Found
$ws.Range("H31").Value = 1794.742
$ws.Range("I31").Value = 1324.3182
$ws.Range("J31").Value = 2944.6667
$ws.Range("K31").Value = 1324.3182
$ws.Range("L31").Value = 2944.6667
$ws.Range("M31").Value = -1029.3182
$ws.Range("N31").Value = -3534.6667

# CRP row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 1794.742
$ws.Range("I34").Value = 1324.3182
$ws.Range("J34").Value = 2944.6667
$ws.Range("K34").Value = 1324.3182
$ws.Range("L34").Value = 2944.6667
$ws.Range("M34").Value = -1122.3182
$ws.Range("N34").Value = -3348.6667

# CRP row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 3079.672
$ws.Range("I58").Value = 1387.4073
$ws.Range("J58").Value = 4423.5293
$ws.Range("K58").Value = 1387.4073
$ws.Range("L58").Value = 4423.5293
$ws.Range("M58").Value = -1184.4073
$ws.Range("N58").Value = -4829.5293

# CRP row 107: Built to Last
$ws.Range("H107").Value = 1263754
$ws.Range("I107").Value = 1603676.6
$ws.Range("J107").Value = 1184
$ws.Range("K107").Value = 1603676.6
$ws.Range("L107").Value = 1184
$ws.Range("M107").Value = -1601756.6
$ws.Range("N107").Value = -5024

# CRP row 132: Hull Lotta Damage
$ws.Range("H132").Value = 2271.4888
$ws.Range("I132").Value = 1957.1082
$ws.Range("K132").Value = 5871.3246
$ws.Range("M132").Value = -3341.3246

# CRP row 134: Wood You Be Quiet
$ws.Range("H134").Value = 1614.1351
$ws.Range("I134").Value = 1216.5312
$ws.Range("K134").Value = 3649.5936
$ws.Range("M134").Value = -1114.5936

# CRP row 136: Turali Quality
$ws.Range("H136").Value = 3079.672
$ws.Range("I136").Value = 1387.4073
$ws.Range("J136").Value = 4423.5293
$ws.Range("K136").Value = 4162.2219
$ws.Range("L136").Value = 13270.5879
$ws.Range("M136").Value = -1612.2219
$ws.Range("N136").Value = -18370.5879

$ws = $wb.Worksheets.Item("CUL")
# CUL row 5: What a Sap
$ws.Range("H5").Value = 11905457
$ws.Range("I5").Value = 508.05264
$ws.Range("J5").Value = 21739980
$ws.Range("K5").Value = 1524.15792
$ws.Range("L5").Value = 65219940
$ws.Range("M5").Value = -1412.15792
$ws.Range("N5").Value = -65220164

# CUL row 64: The Aroma of Faith
$ws.Range("H64").Value = 1001509.7
$ws.Range("J64").Value = 2001842.9
$ws.Range("L64").Value = 6005528.699999999
$ws.Range("N64").Value = -6006068.699999999

# CUL row 67: Soup's On (L)
$ws.Range("H67").Value = 1001509.7
$ws.Range("J67").Value = 2001842.9
$ws.Range("L67").Value = 6005528.699999999
$ws.Range("N67").Value = -6007400.699999999

# CUL row 70: Persona non Gratin
$ws.Range("H70").Value = 23129.889
$ws.Range("I70").Value = 37437.8
$ws.Range("K70").Value = 112313.4
$ws.Range("M70").Value = -111998.4

# CUL row 73: Recipe for Disaster (L)
$ws.Range("H73").Value = 23129.889
$ws.Range("I73").Value = 37437.8
$ws.Range("K73").Value = 112313.4
$ws.Range("M73").Value = -111221.4

# CUL row 80: Saucy for a Suitor
$ws.Range("H80").Value = 10347787
$ws.Range("I80").Value = 24143256
$ws.Range("J80").Value = 1185
$ws.Range("K80").Value = 72429768
$ws.Range("L80").Value = 3555
$ws.Range("M80").Value = -72428832
$ws.Range("N80").Value = -5427

# CUL row 83: Saved by the Sauce (L)
$ws.Range("H83").Value = 10347787
$ws.Range("I83").Value = 24143256
$ws.Range("J83").Value = 1185
$ws.Range("K83").Value = 217289304
$ws.Range("L83").Value = 10665
$ws.Range("M83").Value = -217284624
$ws.Range("N83").Value = -20025

# CUL row 131: The Mountain Steeped
$ws.Range("H131").Value = 3333.7112
$ws.Range("J131").Value = 2700.8708
$ws.Range("L131").Value = 8102.6124
$ws.Range("N131").Value = -18182.6124

# CUL row 135: Not-so-secret Ingredient
$ws.Range("H135").Value = 11905457
$ws.Range("I135").Value = 508.05264
$ws.Range("J135").Value = 21739980
$ws.Range("K135").Value = 4572.47376
$ws.Range("L135").Value = 195659820
$ws.Range("M135").Value = -2037.47376
$ws.Range("N135").Value = -195664890

$ws = $wb.Worksheets.Item("GSM")
# GSM row 107: Whetstones for the Workers
$ws.Range("H107").Value = 611.8095
$ws.Range("I107").Value = 419.5
$ws.Range("J107").Value = 996.4286
$ws.Range("K107").Value = 419.5
$ws.Range("L107").Value = 996.4286
$ws.Range("M107").Value = 1500.5
$ws.Range("N107").Value = -4836.4286

# GSM row 113: Copious Crystal Cannons
$ws.Range("H113").Value = 1477.625
$ws.Range("I113").Value = 1263
$ws.Range("K113").Value = 1263
$ws.Range("M113").Value = 907

# GSM row 132: On Board for Lar
$ws.Range("H132").Value = 2549.75
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 2549.75
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 7649.25
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -12709.25

$ws = $wb.Worksheets.Item("LTW")
# LTW row 55: It's Not a Job, It's a Calling
$ws.Range("H55").Value = 122.833336
$ws.Range("I55").Value = 105.78571
$ws.Range("K55").Value = 105.78571
$ws.Range("M55").Value = 67.21429000000001

# LTW row 61: Spelling Me Softly
$ws.Range("H61").Value = 2200
$ws.Range("I61").Value = 2200
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2200
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1998
$ws.Range("N61").ClearContents()

# LTW row 113: Peace in Rest
$ws.Range("H113").Value = 2200
$ws.Range("I113").Value = 2200
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2200
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -30
$ws.Range("N113").ClearContents()

# LTW row 136: Respect for Br'aax
$ws.Range("H136").Value = 1822
$ws.Range("I136").Value = 1456.8572
$ws.Range("J136").Value = 3100
$ws.Range("K136").Value = 4370.571599999999
$ws.Range("L136").Value = 9300
$ws.Range("M136").Value = -1820.571599999999
$ws.Range("N136").Value = -14400

$ws = $wb.Worksheets.Item("WVR")
# WVR row 113: A Tender Table
$ws.Range("H113").Value = 404.84616
$ws.Range("I113").Value = 404.84616
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1214.53848
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 955.4615200000001
$ws.Range("N113").ClearContents()

# WVR row 131: A Better Bottom Line
$ws.Range("H131").Value = 75312.164
$ws.Range("J131").Value = 75312.164
$ws.Range("L131").Value = 75312.164
$ws.Range("N131").Value = -85392.164

# WVR row 132: Comfy Cabins
$ws.Range("H132").Value = 2616.8948
$ws.Range("I132").Value = 2909.111
$ws.Range("J132").Value = 1899.6364
$ws.Range("K132").Value = 8727.332999999999
$ws.Range("L132").Value = 5698.9092
$ws.Range("M132").Value = -6197.332999999999
$ws.Range("N132").Value = -10758.9092

# WVR row 136: Weaving the Envelope
$ws.Range("H136").Value = 8286.082
$ws.Range("I136").Value = 10835.639
$ws.Range("J136").Value = 1225.7693
$ws.Range("K136").Value = 32506.917
$ws.Range("L136").Value = 3677.3079
$ws.Range("M136").Value = -29956.917
$ws.Range("N136").Value = -8777.3079
